$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# G4 and G6 currently hold the numeric value 10000. Replace them with the
# text string "10 000" (with a non-breaking-looking space), matching the
# author's edit that converts these cells to shared-string text values.
$ws.Range("G4").Value = "10 000"
$ws.Range("G6").Value = "10 000"
